# GFG-Intersection of Two sorted linked list
# Append a new row to the Linked_List question tracker sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 11: same "GFG" source label as the other GFG rows, plus the
# new question title (wrap-text style, matching the row above it).
$ws.Range("A11").Value = "GFG"
$ws.Range("B11").Value = "Intersection of two sorted Linked lists"
$ws.Range("B11").WrapText = $true

# Leave the selection where the author last clicked.
$ws.Range("B6").Select()
